# fixed LBNRIND, docs and delete check sites
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mislabeled header in F1: LBNDIND -> LBNRIND
$ws.Range("F1").Value = "LBNRIND"

# Update the saved selection/active cell to F11
$ws.Range("F11").Select()
